$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2191309.04
$ws.Range("O2").Value = 22996.26

# Row 3
$ws.Range("K3").Value = 835243
$ws.Range("L3").Value = 1848985
$ws.Range("O3").Value = 88144

# Row 4
$ws.Range("K4").Value = 426448.72
$ws.Range("L4").Value = 1158640.7
$ws.Range("O4").Value = 589674.4399999999

# Row 5
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 332002.14
$ws.Range("O5").Value = 2688.86

# Row 6
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 14167357.5
$ws.Range("O6").Value = 7512736.5

# Row 11
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 374093.4
$ws.Range("O11").Value = 1028765.85

# Row 14
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 1272302.27
$ws.Range("O14").Value = 3498887.23
